$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("omnidirectional")
$ws.Range("B14").Value = 101080047.2804286
$ws.Range("C14").Value = 39758439.78462993
$ws.Range("B15").Value = 124696841.3759438
$ws.Range("C15").Value = 47314736.6349548
$ws.Range("B16").Value = 145580277.9066491
$ws.Range("C16").Value = 57394610.27960317
$ws.Range("B17").Value = 168598227.0585145
$ws.Range("C17").Value = 66289703.25172374

$ws = $wb.Worksheets.Item("345 to 15.0")
$ws.Range("B14").Value = 52104034.9396586
$ws.Range("C14").Value = 20927264.4442574
$ws.Range("B15").Value = 62466203.33198615
$ws.Range("C15").Value = 25566838.85353827

$ws = $wb.Worksheets.Item("15.0 to 45.0")
$ws.Range("B14").Value = 40324431.99508546
$ws.Range("C14").Value = 15900105.57303656
$ws.Range("B15").Value = 49359400.70951353
$ws.Range("C15").Value = 19004391.63026109
$ws.Range("B16").Value = 54180446.42290615
$ws.Range("C16").Value = 21513610.97208555
$ws.Range("B17").Value = 58954858.49246357
$ws.Range("C17").Value = 24208919.1374758

$ws = $wb.Worksheets.Item("45.0 to 75.0")
$ws.Range("B15").Value = 49159506.84081858
$ws.Range("C15").Value = 18659792.31229874
$ws.Range("B16").Value = 58994441.38441692
$ws.Range("C16").Value = 21577369.46603119
$ws.Range("B18").Value = 67417486.99980083
$ws.Range("C18").Value = 26326795.12820816

$ws = $wb.Worksheets.Item("75.0 to 105.0")
$ws.Range("B16").Value = 67111096.20276004
$ws.Range("C16").Value = 24424998.42848644
$ws.Range("B19").Value = 88380967.50525707
$ws.Range("C19").Value = 35543069.16547316

$ws = $wb.Worksheets.Item("105.0 to 135.0")
$ws.Range("B16").Value = 60748088.86338413
$ws.Range("C16").Value = 22585184.70832386
$ws.Range("B17").Value = 75251429.41348067
$ws.Range("C17").Value = 26867271.05408737
$ws.Range("B19").Value = 91033072.4913754
$ws.Range("C19").Value = 34084825.68684281
$ws.Range("B20").Value = 89982883.13563851
$ws.Range("C20").Value = 36624429.32719123

$ws = $wb.Worksheets.Item("135.0 to 165.0")
$ws.Range("B18").Value = 64226674.74376989
$ws.Range("C18").Value = 22838308.27617394
$ws.Range("B19").Value = 71784807.48078391
$ws.Range("C19").Value = 24823811.31341133
$ws.Range("B22").Value = 68305157.54598257
$ws.Range("C22").Value = 25229456.21284471
$ws.Range("B23").Value = 59933797.89104021
$ws.Range("C23").Value = 23492625.25433663
$ws.Range("B24").Value = 58593163.11216585
$ws.Range("C24").Value = 24575351.49080865

$ws = $wb.Worksheets.Item("165.0 to 195.0")
$ws.Range("B17").Value = 64993317.05317006
$ws.Range("C17").Value = 23723475.43613141
$ws.Range("B18").Value = 75490241.89233144
$ws.Range("C18").Value = 26598409.0732093
$ws.Range("B20").Value = 90758665.82406816
$ws.Range("C20").Value = 33335861.14041523
$ws.Range("B21").Value = 84414899.8944288
$ws.Range("C21").Value = 33946526.58856174

$ws = $wb.Worksheets.Item("195.0 to 225.0")
$ws.Range("B15").Value = 71929349.37207317
$ws.Range("C15").Value = 27155221.94164393
$ws.Range("B17").Value = 96837931.10622483
$ws.Range("C17").Value = 37466850.59425086
$ws.Range("B18").Value = 104050069.5000892
$ws.Range("C18").Value = 42386095.63056424

$ws = $wb.Worksheets.Item("225.0 to 255.0")
$ws.Range("B14").Value = 67476434.57885192
$ws.Range("C14").Value = 26600521.16821748
$ws.Range("B16").Value = 99711716.6615762
$ws.Range("C16").Value = 39683389.52205569
$ws.Range("B17").Value = 107771815.1815404
$ws.Range("C17").Value = 45578651.36204857

$ws = $wb.Worksheets.Item("255.0 to 285.0")
$ws.Range("B14").Value = 64761874.18430185
$ws.Range("C14").Value = 25729585.83296733
$ws.Range("B15").Value = 75987788.15837739
$ws.Range("C15").Value = 31264075.37124481
$ws.Range("B16").Value = 90820153.21959911
$ws.Range("C16").Value = 37627259.25389922

$ws = $wb.Worksheets.Item("285.0 to 315.0")
$ws.Range("B14").Value = 62354363.96756563
$ws.Range("C14").Value = 24828532.35635889
$ws.Range("B15").Value = 73578207.32170214
$ws.Range("C15").Value = 30241877.50997814
$ws.Range("B16").Value = 85625636.96112494
$ws.Range("C16").Value = 36149867.3946562

$ws = $wb.Worksheets.Item("315.0 to 345.0")
$ws.Range("B13").Value = 52316963.5169438
$ws.Range("C13").Value = 20655344.31762617
$ws.Range("B14").Value = 63499507.22019055
$ws.Range("C14").Value = 25708466.36669232
$ws.Range("B15").Value = 73382773.42880546
$ws.Range("C15").Value = 30458145.75722643

